$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74-129 down to 75-130.
$ws.Rows(74).Insert()

# Populate the newly inserted row 74 with the new data record.
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44566
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100103
$ws.Range("H74").Value = "Frutos de hueso (carozo)"
$ws.Range("I74").Value = 100103004
$ws.Range("J74").Value = "Durazno"
$ws.Range("K74").Value = "Kurakata"
$ws.Range("L74").Value = "Segunda"
$ws.Range("M74").Value = 220
$ws.Range("N74").Value = 12000
$ws.Range("O74").Value = 13000
$ws.Range("P74").Value = 12455
$ws.Range("Q74").Value = "$/caja 16 kilos empedrada"
$ws.Range("R74").Value = "Provincia de Curicó"
$ws.Range("S74").Value = 778
$ws.Range("T74").Value = 16
